# Slide 6 ("Analyse des résultats") contains a single run whose text reads:
# " (dont le but est la synchronisation) sans laissant le temps à la cible
#   de répondre, le tout en laissant le port ouvert et en envoyant une autre
#   demande sur un autre port, et ainsi de suite."
# The author's edit changes "le port ouvert" to "le port semi-ouvert" inside
# that run, leaving everything else (run formatting, other text, other
# slides) untouched.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

$fullText = $tr.Text

# Locate the run that needs editing using the stable text immediately
# before/after the word that changes, so we replace the whole run's text
# in one shot (this keeps it as a single <a:r> run in the XML, matching
# how the original author's edit looks).
$runStartMarker = " (dont le but est la synchronisation)"
$runEndMarker   = "ainsi de suite."

$startIdx = $fullText.IndexOf($runStartMarker)
if ($startIdx -lt 0) {
    throw "Could not locate target run start marker on slide 6"
}

$endIdx = $fullText.IndexOf($runEndMarker, $startIdx)
if ($endIdx -lt 0) {
    throw "Could not locate target run end marker on slide 6"
}

$runLength = ($endIdx + $runEndMarker.Length) - $startIdx
$runRange  = $tr.Characters($startIdx + 1, $runLength)

$oldRunText = $runRange.Text
$newRunText = $oldRunText.Replace("le port ouvert", "le port semi-ouvert")

if ($newRunText -ne $oldRunText) {
    $runRange.Text = $newRunText
}
